# Populate counter PLL calculator w/ 640x480 @ 60 Hz values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Counter PLL Calc")

# Desired counter output frequency (pixel clock for 640x480 @ 60 Hz)
$ws.Range("C6").Value = 25175000
